$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A from 45 to 47 (the stored OOXML width differs from the
# COM ColumnWidth property by Excel's internal MDW-based padding; 46.17
# is calibrated to round-trip to a stored width of exactly 47)
$ws.Columns.Item(1).ColumnWidth = 46.17

# Append new Q&A rows (22-26)
$ws.Range("A22").Value = "What are the formats for loading a text file?"
$ws.Range("B22").Value = "Yes, GEO can load several different types of ASCII files, including tab-delimited and space-delimited files. In addition to these files, also csv (comma-separated value) files can be loaded."

$ws.Range("A23").Value = "What is an ODT?"
$ws.Range("B23").Value = "An ODT (OpenDocument Template) is a file format used in the GEO application that remains unchanged after generating an ODF (OpenDocument File) from a template, allowing users to reuse it for subsequent wells."

$ws.Range("A24").Value = "What is an ODT?"
$ws.Range("B24").Value = "An ODT (OpenDocument Template) is a file format used in the GEO application that remains unchanged after generating an ODF (OpenDocument File) from a template, allowing users to reuse it for subsequent wells."

$ws.Range("A25").Value = "What is the limit on the number of curves?"
$ws.Range("B25").Value = "The limit on the number of curves is 450, with an additional note that the infinite wrapping is limited by the `"Maximum number of times to Wrap a Curve`" in the Global Settings."

$ws.Range("A26").Value = "Where can you view the mouse pointer depth?"
$ws.Range("B26").Value = "To view the mouse pointer depth, follow these steps:
1. Go to the Depth tab and select Index Converter or press + on your keyboard.
2. In the Convert Index dialog box, enter the desired depth or time of the section of log you wish to view.
3. Click `"Go`" to instantly display the relevant depth section.
Alternatively, you can also access this dialog box by going to Tools and selecting Index Converter or using the Mouse & Keyboard Modestoolbar."

# The multi-line text in B26 causes the runtime to auto-expand the row
# height; AutoFit restores it back to the sheet's default (no explicit
# row height / customHeight attribute), matching the source row style.
$ws.Rows.Item(26).AutoFit()
